# "recording the war data"
#
# Record_Cooldown sheet: add two new config rows ("Ref" and "Force", both
# with a value of 0) right before the existing "Upload" row. That pushes
# "Upload", the "SkillID"/"Time" header row, the "string"/"int" type row and
# the trailing "Desc" row down by two rows (8->10, 9->11, 10->12, 11->13),
# and the previously-empty B11 cell (now B13) gets filled in with "Desc" too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record_Cooldown")

# Shift the four trailing rows down by two (bottom-up so we never overwrite
# a row before it has been copied). Copying full A:B cells this way carries
# over both the value and the existing style, so it lines up with the
# formatting already used on this sheet instead of inventing new styles.
$ws.Range("A11:B11").Copy($ws.Range("A13:B13"))
$ws.Range("A10:B10").Copy($ws.Range("A12:B12"))
$ws.Range("A9:B9").Copy($ws.Range("A11:B11"))
$ws.Range("A8:B8").Copy($ws.Range("A10:B10"))
$excel.CutCopyMode = $false

# Fill the two freshly-freed rows (8 and 9) with the new "Ref" / "Force"
# entries, reusing the same formatting as the other whole-number rows above
# them (e.g. row 7, "Cache") instead of leaving the default style behind.
$ws.Range("A7:B7").Copy($ws.Range("A8:B8"))
$ws.Range("A7:B7").Copy($ws.Range("A9:B9"))
$excel.CutCopyMode = $false

$ws.Cells.Item(8, 1).Value = "Ref"
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(9, 1).Value = "Force"
$ws.Cells.Item(9, 2).Value = 0

# The row that used to be row 11 (A11="Desc", B11 blank) is now row 13;
# fill in the previously-blank B13 with "Desc" as well.
$ws.Cells.Item(13, 2).Value = "Desc"

# Match the authored selection/active cell.
$ws.Range("A10").Select()
